# Update pl_mw.xlsx results for the 380 kV line case (Case_0_167).
# Columns C..O (3..15) of rows 2..25 are rewritten with newly computed
# per-unit line power-flow results; columns A, B and H are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    @(0.1710253105893287,0.1249055367452456,0.135974826705981,1.901870967445312,0.002505621797156709,0,1.340883869616349,0.1798759958227016,1.292249610345806,0.1937768954559367,0.4476476930700457,1.638662744124019,4.99407545427303),
    @(0.1693807721009932,0.1233959450520814,0.136137818715687,1.908938088328789,0.002508442882714232,0,1.346195998359519,0.1809549444737062,1.19432589802895,0.1946971206182013,0.4287578226186852,1.648283506349223,5.020225116108264),
    @(0.1684304603477642,0.122505840444731,0.1362800952807675,1.914188079044067,0.002510268882654455,0,1.350107772849881,0.1816726150959358,1.1343845162354,0.1953157865117241,0.4172813853441966,1.65473330055768,5.038773278176023),
    @(0.1680582148325556,0.1221524244316541,0.136348712771909,1.916556645055138,0.00251103665955744,0,1.351865410487818,0.1819789692126577,1.110006011299618,0.1955814121080408,0.4126357227814665,1.657498307186451,5.04695861395831),
    @(0.1679973125492324,0.1220943038958495,0.1363607500079222,1.916963787904308,0.002511165579843729,0,1.352167145831771,0.1820306790179167,1.105960936481438,0.1956263359983552,0.4118662028032247,1.657965695110008,5.048355644475009),
    @(0.1684253792326587,0.1225010363860264,0.1362809775682869,1.914219094345135,0.002510279141196926,0,1.350130814616918,0.1816766903989073,1.134055542093392,0.1953193140791143,0.4172186059519447,1.654770036765981,5.038881130377547),
    @(0.1704459760801171,0.1243774273938598,0.1360222842269447,1.904118763811724,0.002506575072666332,0,1.342580610214078,0.1802365737129534,1.25844832641485,0.1940830702650906,0.4411093714341305,1.641867495176115,5.002574684002269),
    @(0.1748773874356857,0.1283466849387764,0.135848691187439,1.891533780485389,0.002500052813822208,0,1.332930253367493,0.1778495940433693,1.503777573720299,0.1920833564991717,0.4889124489919396,1.620861685957927,4.951149087454951),
    @(0.1784160585480805,0.1314368428752601,0.1359231915007122,1.88668589378004,0.002495708429966687,0,1.328980828394378,0.176361214934186,1.684801206155953,0.1908715533690852,0.5245984667571335,1.608035705056878,4.925420142275868),
    @(0.1800867157179198,0.1328798719842155,0.1360006675226799,1.885434768472052,0.002493828283045103,0,1.327865731415322,0.1757414927627643,1.767309174221282,0.1903758619921661,0.5409524932380876,1.602764510772218,4.916333125652102),
    @(0.1807280461847824,0.133431620589036,0.1360362475395149,1.885098121193067,0.002493130072638895,0,1.327541423106922,0.1755150485422341,1.798574155217466,0.1901961224050979,0.5471622840003576,1.600849278343787,4.91326844432416),
    @(0.1805895388405219,0.133312556660897,0.1360283075659154,1.885164526835382,0.002493279833610582,0,1.327606912880647,0.1755634515228621,1.791839778496353,0.1902344785128527,0.5458241505711072,1.601258164429865,4.91391173854376),
    @(0.1801393046143858,0.1329251586946327,0.1360034697711114,1.885404324568896,0.002493770565439202,0,1.327837087483232,0.1757226981856235,1.76988095304506,0.190360915143696,0.5414630406463559,1.602605323918667,4.916073449232016),
    @(0.1798646526747945,0.1326885553476416,0.1359890679655749,1.885569062987727,0.002494072942812494,0,1.327990831290663,0.1758213128382806,1.756433216896596,0.1904393982451609,0.5387939201798204,1.603441022894074,4.917446575231111),
    @(0.1783080967632173,0.1313432830688441,0.1359190030863218,1.886786851364448,0.002495833230997356,0,1.329067412477393,0.1764028679075551,1.679412160008042,0.1909050641675023,0.5235320770454877,1.608391513161592,4.926066667092186),
    @(0.177368748990844,0.1305275136927264,0.1358871653028686,1.887778257724705,0.00249693768784138,0,1.329902379719329,0.1767743104144444,1.632201698308791,0.1912049506592872,0.5141999375905186,1.611572663339707,4.932025160670236),
    @(0.1768341970702778,0.1300618191676719,0.1358729570713813,1.888438307665567,0.002497581994798659,0,1.330446781095908,0.1769933531988706,1.605062608213302,0.1913826686892435,0.508843683982775,1.613455416872888,4.935698671290169),
    @(0.1766541942294566,0.1299047484120663,0.1358688520712477,1.888677218412951,0.00249780170268739,0,1.33064212534498,0.1770684450622646,1.595876433832018,0.1914437401070153,0.5070321100740003,1.614101999616224,4.936984767781951),
    @(0.1774681509975125,0.1306139904316126,0.1358901299273718,1.887663425497209,0.00249681918017599,0,1.32980685705099,0.176734211075221,1.637225781723544,0.1911724859754145,0.5151921885244519,1.611228536219187,4.931365373138021),
    @(0.1802713140506,0.1330388033675263,0.1360105960471536,1.885330169270304,0.002493626052641529,0,1.327766821530098,0.1756757003139988,1.776330237593129,0.1903235615667782,0.542743549682406,1.602207437426813,4.915428287515283),
    @(0.1821539568700246,0.1346544519412447,0.1361256960153163,1.884604490897971,0.002491619336821728,0,1.327004473931439,0.1750318724377671,1.867364581801041,0.1898151737316596,0.5608481280852828,1.596782837666922,4.907206291412109),
    @(0.1811445463101222,0.1337893426254766,0.1360609458571851,1.884918697173561,0.002492683043444125,0,1.327359128073439,0.1753711113454415,1.818767350863197,0.1900822684853054,0.5511765297183473,1.599634985382096,4.911393783879475),
    @(0.1774231942401912,0.1305748840079843,0.1358887768611545,1.887715060554385,0.002496872728650524,0,1.329849842343663,0.1767523228671202,1.634954384639059,0.1911871467194892,0.51474356406748,1.611383948280583,4.931662890903453),
    @(0.1736286731594845,0.1272421589495636,0.1358600720176124,1.894165657787283,0.002501738348091969,0,1.334989228278971,0.1784486590085699,1.437267425543098,0.1925790291399174,0.475880053976816,1.62608565517457,4.962944394620791)
)

$firstRow = 2
$firstCol = 3   # column C

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $r = $firstRow + $i
    $rowVals = $newValues[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $c = $firstCol + $j
        $ws.Cells.Item($r, $c).Value2 = $rowVals[$j]
    }
}

Write-Host "case with 380 kV done"
